# Auto-generated edit script: updates Leve profit-calculation cells
# (H/I/J/K/L/M/N columns) across 8 worksheets to reflect refreshed
# market-price data, per the scheduled-runner commit.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 86
$ws.Range("H86").Value = 7806.1333
$ws.Range("I86").Value = 1295.4546
$ws.Range("J86").Value = 25710.5
$ws.Range("K86").Value = 1295.4546
$ws.Range("L86").Value = 25710.5
$ws.Range("M86").Value = -172.4546
$ws.Range("N86").Value = -27956.5
# Row 89
$ws.Range("H89").Value = 7806.1333
$ws.Range("I89").Value = 1295.4546
$ws.Range("J89").Value = 25710.5
$ws.Range("K89").Value = 6477.273
$ws.Range("L89").Value = 128552.5
$ws.Range("M89").Value = -861.2730000000001
$ws.Range("N89").Value = -139784.5
# Row 106
$ws.Range("H106").Value = 11497405
$ws.Range("I106").Value = 41668556
$ws.Range("J106").Value = 3634.1904
$ws.Range("K106").Value = 41668556
$ws.Range("L106").Value = 3634.1904
$ws.Range("M106").Value = -41667925
$ws.Range("N106").Value = -4896.190399999999
# Row 107
$ws.Range("H107").Value = 882.5238000000001
$ws.Range("I107").Value = 826.5333000000001
$ws.Range("J107").Value = 1022.5
$ws.Range("K107").Value = 826.5333000000001
$ws.Range("L107").Value = 1022.5
$ws.Range("M107").Value = 1093.4667
$ws.Range("N107").Value = -4862.5
# Row 129
$ws.Range("H129").Value = 724.0333000000001
$ws.Range("J129").Value = 805.36
$ws.Range("L129").Value = 2416.08
$ws.Range("N129").Value = -12416.08
# Row 132
$ws.Range("H132").Value = 3095.3635
$ws.Range("I132").Value = 3408.6072
$ws.Range("J132").Value = 1341.2
$ws.Range("K132").Value = 10225.8216
$ws.Range("L132").Value = 4023.6
$ws.Range("M132").Value = -7695.821599999999
$ws.Range("N132").Value = -9083.6
# Row 137
$ws.Range("H137").Value = 1858.8235
$ws.Range("I137").Value = 2115
$ws.Range("J137").Value = 1492.8572
$ws.Range("K137").Value = 6345
$ws.Range("L137").Value = 4478.571599999999
$ws.Range("M137").Value = -3795
$ws.Range("N137").Value = -9578.571599999999

$ws = $wb.Worksheets.Item("ARM")
# Row 7
$ws.Range("H7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").ClearContents()
# Row 32
$ws.Range("H32").Value = 5926.6895
$ws.Range("I32").Value = 4632.71
$ws.Range("J32").Value = 10886.944
$ws.Range("K32").Value = 4632.71
$ws.Range("L32").Value = 10886.944
$ws.Range("M32").Value = -4345.71
$ws.Range("N32").Value = -11460.944
# Row 61
$ws.Range("H61").Value = 3576
$ws.Range("I61").Value = 3336.3635
$ws.Range("J61").Value = 5333.3335
$ws.Range("K61").Value = 3336.3635
$ws.Range("L61").Value = 5333.3335
$ws.Range("M61").Value = -3124.3635
$ws.Range("N61").Value = -5757.3335
# Row 110
$ws.Range("H110").Value = 836.875
$ws.Range("I110").Value = 738.2143
$ws.Range("K110").Value = 738.2143
$ws.Range("M110").Value = 1306.7857
# Row 136
$ws.Range("H136").Value = 3576
$ws.Range("I136").Value = 3336.3635
$ws.Range("J136").Value = 5333.3335
$ws.Range("K136").Value = 10009.0905
$ws.Range("L136").Value = 16000.0005
$ws.Range("M136").Value = -7459.0905
$ws.Range("N136").Value = -21100.0005

$ws = $wb.Worksheets.Item("BSM")
# Row 51
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()
# Row 94
$ws.Range("H94").Value = 576.74194
$ws.Range("I94").Value = 433
$ws.Range("K94").Value = 433
$ws.Range("M94").Value = 18
# Row 105
$ws.Range("H105").Value = 1725981
$ws.Range("I105").Value = 1695
$ws.Range("J105").Value = 2633500
$ws.Range("K105").Value = 1695
$ws.Range("L105").Value = 2633500
$ws.Range("M105").Value = 52
$ws.Range("N105").Value = -2636994
# Row 107
$ws.Range("H107").Value = 1892
$ws.Range("I107").Value = 1980.3334
$ws.Range("J107").Value = 1839
$ws.Range("K107").Value = 1980.3334
$ws.Range("L107").Value = 1839
$ws.Range("M107").Value = -60.33339999999998
$ws.Range("N107").Value = -5679

$ws = $wb.Worksheets.Item("CRP")
# Row 105
$ws.Range("H105").Value = 1288.6428
$ws.Range("J105").Value = 2015.25
$ws.Range("L105").Value = 2015.25
$ws.Range("N105").Value = -5509.25
# Row 134
$ws.Range("H134").Value = 1272.1052
$ws.Range("I134").Value = 1122.8572
$ws.Range("J134").Value = 1690
$ws.Range("K134").Value = 3368.5716
$ws.Range("L134").Value = 5070
$ws.Range("M134").Value = -833.5715999999998
$ws.Range("N134").Value = -10140

$ws = $wb.Worksheets.Item("CUL")
# Row 40
$ws.Range("H40").Value = 106.44444
$ws.Range("I40").Value = 96.85714
$ws.Range("J40").Value = 140
$ws.Range("K40").Value = 387.42856
$ws.Range("L40").Value = 560
$ws.Range("M40").Value = -318.42856
$ws.Range("N40").Value = -698
# Row 131
$ws.Range("H131").Value = 714.6288500000001
$ws.Range("J131").Value = 742.6853599999999
$ws.Range("L131").Value = 2228.05608
$ws.Range("N131").Value = -12308.05608
# Row 133
$ws.Range("H133").Value = 5599.364
$ws.Range("I133").Value = 1315
$ws.Range("J133").Value = 6551.4443
$ws.Range("K133").Value = 3945
$ws.Range("L133").Value = 19654.3329
$ws.Range("M133").Value = 1115
$ws.Range("N133").Value = -29774.3329

$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 47.88889
$ws.Range("I2").Value = 40.857143
$ws.Range("K2").Value = 40.857143
$ws.Range("M2").Value = 72.14285699999999
# Row 102
$ws.Range("H102").Value = 1694.2307
$ws.Range("I102").Value = 1355.3043
$ws.Range("K102").Value = 1355.3043
$ws.Range("M102").Value = 266.6957
# Row 113
$ws.Range("H113").Value = 4920.433
$ws.Range("I113").Value = 6083.75
$ws.Range("J113").Value = 2593.8
$ws.Range("K113").Value = 6083.75
$ws.Range("L113").Value = 2593.8
$ws.Range("M113").Value = -3913.75
$ws.Range("N113").Value = -6933.8
# Row 126
$ws.Range("H126").Value = 2723.1162
$ws.Range("I126").Value = 2160.4783
$ws.Range("K126").Value = 6481.4349
$ws.Range("M126").Value = -4011.4349
# Row 132
$ws.Range("H132").Value = 28803.334
$ws.Range("I132").Value = 5304.8125
$ws.Range("J132").Value = 103998.6
$ws.Range("K132").Value = 15914.4375
$ws.Range("L132").Value = 311995.8
$ws.Range("M132").Value = -13384.4375
$ws.Range("N132").Value = -317055.8

$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 2810.7407
$ws.Range("I40").Value = 2443.0952
$ws.Range("J40").Value = 4097.5
$ws.Range("K40").Value = 2443.0952
$ws.Range("L40").Value = 4097.5
$ws.Range("M40").Value = -2307.0952
$ws.Range("N40").Value = -4369.5
# Row 122
$ws.Range("H122").Value = 703522.8
$ws.Range("I122").Value = 819934.9399999999
$ws.Range("J122").Value = 5050
$ws.Range("K122").Value = 2459804.82
$ws.Range("L122").Value = 15150
$ws.Range("M122").Value = -2457354.82
$ws.Range("N122").Value = -20050

$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 1737.4
$ws.Range("I122").Value = 1719.3334
$ws.Range("J122").Value = 1900
$ws.Range("K122").Value = 5158.0002
$ws.Range("L122").Value = 5700
$ws.Range("M122").Value = -2708.0002
$ws.Range("N122").Value = -10600
# Row 126
$ws.Range("H126").Value = 1369.5135
$ws.Range("I126").Value = 1012.0323
$ws.Range("K126").Value = 3036.0969
$ws.Range("M126").Value = -566.0969
